$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Dades_Meteo")

# Widen column L (column 12) from 25 to 26 raw OOXML width units.
# Excel's ColumnWidth (character units) round-trips to the stored OOXML
# width with a fixed +5/6 offset on this sheet's font, so back it out.
$ws.Columns.Item(12).ColumnWidth = 26 - 5/6

# Plain text field refreshes for row 2 (these don't look like numbers,
# so Excel keeps them as text automatically).
$ws.Range("E2").Value = "2026-02-20 10:47:37"
$ws.Range("J2").Value = "1021.0 hPa"
$ws.Range("K2").Value = "3.8 MJ/m2"
$ws.Range("L2").Value = "8.6 km/h - 160º 10:02 TU"
$ws.Range("M2").Value = "14.4 °C 10:29 TU"
$ws.Range("O2").Value = "3.8 °C"

# H2 ("89%") looks like a percentage, so a plain .Value assignment would
# have Excel auto-convert it into a numeric percent cell. Use a
# quote-prefixed entry to force text, then restore the original
# (non quote-prefixed) cell formatting by pasting formats from a
# neighboring plain-text cell on the same row.
$ws.Range("H2").Value = "'89%"
$ws.Range("I2").Copy()
$ws.Range("H2").PasteSpecial(-4122)
$excel.CutCopyMode = $false
